$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.295.93"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -4.34%  "

$ws.Range("D3").Value = "'2.502.70"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -5.44%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'575.61"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.30%  "

$ws.Range("D6").Value = "'166.08"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.93%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "'0.515"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.82%  "

$ws.Range("D9").Value = "'2.501.14"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.48%  "

$ws.Range("E10").Value = "  -9.36%  "

$ws.Range("E11").Value = "  -1.30%  "

$ws.Range("E12").Value = "  -4.12%  "

$ws.Range("D13").Value = "'4.83"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.40%  "

$ws.Range("D14").Value = "'2.960.20"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.40%  "

$ws.Range("D15").Value = "'69.254.00"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.17%  "

$ws.Range("E16").Value = "  -7.11%  "

$ws.Range("D17").Value = "'24.70"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.83%  "

$ws.Range("D18").Value = "'2.502.15"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.06%  "

$ws.Range("D19").Value = "'11.38"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -6.30%  "

$ws.Range("E20").Value = "  -2.98%  "

$ws.Range("D21").Value = "'346.96"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.68%  "

$ws.Range("E22").Value = "  -5.86%  "

$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").Value = "'1.94"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -6.10%  "

$ws.Range("D25").Value = "'68.19"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.95%  "

$ws.Range("E26").Value = "  -7.20%  "

$ws.Range("D27").Value = "'8.87"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -8.13%  "

$ws.Range("E28").Value = "  -5.10%  "

$ws.Range("D29").Value = "'1.01"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.48%  "

$ws.Range("D30").Value = "'0.0₃0893"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.73%  "

$ws.Range("D31").Value = "'7.82"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.42%  "

$ws.Range("D32").Value = "'465.74"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -6.49%  "

$ws.Range("E33").Value = "  -2.28%  "

$ws.Range("D34").Value = "'1.74"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.42%  "

$ws.Range("E36").Value = "  +0.75%  "

$ws.Range("D37").Value = "'152.85"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.94%  "

$ws.Range("D38").Value = "'18.94"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.32%  "

$ws.Range("E39").Value = "  -5.00%  "

$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("D41").Value = "'4.72"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.92%  "

$ws.Range("D42").Value = "'0.315"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.99%  "

$ws.Range("E43").Value = "  -9.19%  "

$ws.Range("D44").Value = "'1.16"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -15.00%  "

$ws.Range("D45").Value = "'2.30"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -11.04%  "

$ws.Range("D46").Value = "'38.06"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.53%  "

$ws.Range("D47").Value = "'142.53"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.73%  "

$ws.Range("E48").Value = "  -4.48%  "

$ws.Range("E49").Value = "  -4.66%  "

$ws.Range("D50").Value = "'1.59"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.37%  "

$ws.Range("E51").Value = "  -2.25%  "
